$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20
$data[0,0] = "ECs"
$data[0,1] = "Trf"
$data[0,2] = "Tfr2"
$data[0,3] = "FAPs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 0.7939700000000002
$data[0,7] = 2.38191
$data[0,8] = 0.008990908376056485
$data[0,9] = 0.008990908376056485
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.4093263333333334
$data[0,13] = 1.227979
$data[0,14] = 0.6924542622754842
$data[0,15] = 0.6924542622754842
$data[0,16] = 0.3249928288766668
$data[0,17] = 2.924935459890001
$data[0,18] = 0.006225792826728666
$data[0,19] = 0.006225792826728666
$data[1,0] = "ECs"
$data[1,1] = "Trf"
$data[1,2] = "Tfr2"
$data[1,3] = "Resolving-Mac"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.7939700000000002
$data[1,7] = 2.38191
$data[1,8] = 0.008990908376056485
$data[1,9] = 0.008990908376056485
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.1817976666666667
$data[1,13] = 0.545393
$data[1,14] = 0.3075457377245157
$data[1,15] = 0.3075457377245158
$data[1,16] = 0.1443418934033334
$data[1,17] = 1.29907704063
$data[1,18] = 0.002765115549327819
$data[1,19] = 0.00276511554932782
$data[2,0] = "FAPs"
$data[2,1] = "Trf"
$data[2,2] = "Tfr2"
$data[2,3] = "FAPs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 4.248509666666666
$data[2,7] = 12.745529
$data[2,8] = 0.04811008117156854
$data[2,9] = 0.04811008117156854
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.4093263333333334
$data[2,13] = 1.227979
$data[2,14] = 0.6924542622754842
$data[2,15] = 0.6924542622754842
$data[2,16] = 1.739026883987889
$data[2,17] = 15.651241955891
$data[2,18] = 0.03331403076567215
$data[2,19] = 0.03331403076567215
$data[3,0] = "FAPs"
$data[3,1] = "Trf"
$data[3,2] = "Tfr2"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 4.248509666666666
$data[3,7] = 12.745529
$data[3,8] = 0.04811008117156854
$data[3,9] = 0.04811008117156854
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.1817976666666667
$data[3,13] = 0.545393
$data[3,14] = 0.3075457377245157
$data[3,15] = 0.3075457377245158
$data[3,16] = 0.7723691442107776
$data[3,17] = 6.951322297897
$data[3,18] = 0.01479605040589638
$data[3,19] = 0.01479605040589638
$data[4,0] = "Inflammatory-Mac"
$data[4,1] = "Trf"
$data[4,2] = "Tfr2"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 23.67666533333333
$data[4,7] = 71.029996
$data[4,8] = 0.2681143225342933
$data[4,9] = 0.2681143225342933
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.4093263333333334
$data[4,13] = 1.227979
$data[4,14] = 0.6924542622754842
$data[4,15] = 0.6924542622754842
$data[4,16] = 9.691482606453778
$data[4,17] = 87.223343458084
$data[4,18] = 0.1856569054159753
$data[4,19] = 0.1856569054159753
$data[5,0] = "Inflammatory-Mac"
$data[5,1] = "Trf"
$data[5,2] = "Tfr2"
$data[5,3] = "Resolving-Mac"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 23.67666533333333
$data[5,7] = 71.029996
$data[5,8] = 0.2681143225342933
$data[5,9] = 0.2681143225342933
$data[5,10] = 1
$data[5,11] = 0.3333333333333333
$data[5,12] = 0.1817976666666667
$data[5,13] = 0.545393
$data[5,14] = 0.3075457377245157
$data[5,15] = 0.3075457377245158
$data[5,16] = 4.304362512047556
$data[5,17] = 38.739262608428
$data[5,18] = 0.08245741711831797
$data[5,19] = 0.08245741711831799
$data[6,0] = "MuSCs"
$data[6,1] = "Trf"
$data[6,2] = "Tfr2"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.2020363333333333
$data[6,7] = 0.606109
$data[6,8] = 0.002287857427402051
$data[6,9] = 0.002287857427402051
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.4093263333333334
$data[6,13] = 1.227979
$data[6,14] = 0.6924542622754842
$data[6,15] = 0.6924542622754842
$data[6,16] = 0.08269879152344446
$data[6,17] = 0.7442891237110001
$data[6,18] = 0.001584236627083175
$data[6,19] = 0.001584236627083174
$data[7,0] = "MuSCs"
$data[7,1] = "Trf"
$data[7,2] = "Tfr2"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.2020363333333333
$data[7,7] = 0.606109
$data[7,8] = 0.002287857427402051
$data[7,9] = 0.002287857427402051
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.1817976666666667
$data[7,13] = 0.545393
$data[7,14] = 0.3075457377245157
$data[7,15] = 0.3075457377245158
$data[7,16] = 0.03672973398188889
$data[7,17] = 0.330567605837
$data[7,18] = 0.0007036208003188765
$data[7,19] = 0.0007036208003188765
$data[8,0] = "Neutrophils"
$data[8,1] = "Trf"
$data[8,2] = "Tfr2"
$data[8,3] = "FAPs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 41.77897400000001
$data[8,7] = 125.336922
$data[8,8] = 0.4731046856677785
$data[8,9] = 0.4731046856677785
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.4093263333333334
$data[8,13] = 1.227979
$data[8,14] = 0.6924542622754842
$data[8,15] = 0.6924542622754842
$data[8,16] = 17.10123423784867
$data[8,17] = 153.911108140638
$data[8,18] = 0.3276033560931564
$data[8,19] = 0.3276033560931564
$data[9,0] = "Neutrophils"
$data[9,1] = "Trf"
$data[9,2] = "Tfr2"
$data[9,3] = "Resolving-Mac"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 41.77897400000001
$data[9,7] = 125.336922
$data[9,8] = 0.4731046856677785
$data[9,9] = 0.4731046856677785
$data[9,10] = 1
$data[9,11] = 0.3333333333333333
$data[9,12] = 0.1817976666666667
$data[9,13] = 0.545393
$data[9,14] = 0.3075457377245157
$data[9,15] = 0.3075457377245158
$data[9,16] = 7.595319988927334
$data[9,17] = 68.35787990034601
$data[9,18] = 0.1455013295746221
$data[9,19] = 0.1455013295746221
$data[10,0] = "Resolving-Mac"
$data[10,1] = "Trf"
$data[10,2] = "Tfr2"
$data[10,3] = "FAPs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 17.60794066666666
$data[10,7] = 52.82382199999999
$data[10,8] = 0.1993921448229012
$data[10,9] = 0.1993921448229012
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.4093263333333334
$data[10,13] = 1.227979
$data[10,14] = 0.6924542622754842
$data[10,15] = 0.6924542622754842
$data[10,16] = 7.207393790637556
$data[10,17] = 64.86654411573799
$data[10,18] = 0.1380699405468686
$data[10,19] = 0.1380699405468686
$data[11,0] = "Resolving-Mac"
$data[11,1] = "Trf"
$data[11,2] = "Tfr2"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 17.60794066666666
$data[11,7] = 52.82382199999999
$data[11,8] = 0.1993921448229012
$data[11,9] = 0.1993921448229012
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.1817976666666667
$data[11,13] = 0.545393
$data[11,14] = 0.3075457377245157
$data[11,15] = 0.3075457377245158
$data[11,16] = 3.201082528005111
$data[11,17] = 28.809742752046
$data[11,18] = 0.06132220427603264
$data[11,19] = 0.06132220427603265

$ws.Range("A2:T13").Value = $data
